$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "Data Type" column (E), which
# becomes the new "Scaling Factor" column. This shifts Data Type -> F and
# Topic Name Suffix -> G.
$ws.Columns("E:E").Insert()

# Scaling factor formulas for each metric row.
$ws.Range("E2").Formula = "=POWER(2,-16)"
$ws.Range("E3").Formula = "=POWER(2,-16)"

# Update the endianness / data type values from INT16 to INT32 (now in
# column F after the insert).
$ws.Range("F2").Value = "INT32"
$ws.Range("F3").Value = "INT32"

# New header for the inserted column.
$ws.Range("E1").Value = "Scaling Factor"
